# Executive_Performance_Summary.xlsx — pivot rebuild:
#   - re-sort entities alphabetically (CryptoFlow, FinShield Re, Nexus
#     Strategic, Omni-Retail, Terra-Grid), dropping the BioGrowth row
#     (negative-NOPAT entity excluded from the final view)
#   - swap the NOPAT/Assets columns so C=assets, D=nopat and lower-case
#     every header label
#   - recompute NOPAT as 75% of revenue and round the sample ROIC figure
#   - Terra-Grid's ROIC_% cell becomes a real 0 instead of a blank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BioGrowth row (old row 4) no longer belongs in the summary - remove
# the whole row and let the rows below shift up.
$ws.Rows("4").Delete()

# --- Header row (now lower-case, NOPAT/Assets swapped) -------------------
$ws.Range("A1").Value = "entity"
$ws.Range("B1").Value = "revenue"
$ws.Range("C1").Value = "assets"
$ws.Range("D1").Value = "nopat"
$ws.Range("E1").Value = "roic_%"

# --- Row 2: CryptoFlow ----------------------------------------------------
$ws.Range("A2").Value = "CryptoFlow"
$ws.Range("B2").Value = 423000
$ws.Range("C2").Value = 189000
$ws.Range("D2").Value = 317250
$ws.Range("E2").Value = 167.86

# --- Row 3: FinShield Re --------------------------------------------------
# This physical row used to be Terra-Grid (a numeric 0 in ROIC_%); the new
# occupant has no ROIC figure, so blank the leftover number back out.
$ws.Range("A3").Value = "FinShield Re"
$ws.Range("B3").Value = 5270500
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 3952875
$ws.Range("E3").ClearContents()

# --- Row 4: Nexus Strategic -----------------------------------------------
$ws.Range("A4").Value = "Nexus Strategic"
$ws.Range("B4").Value = 869200
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 651900

# --- Row 5: Omni-Retail ----------------------------------------------------
$ws.Range("A5").Value = "Omni-Retail"
$ws.Range("B5").Value = 9483000
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 7112250

# --- Row 6: Terra-Grid ------------------------------------------------------
$ws.Range("A6").Value = "Terra-Grid"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 2100000
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
